$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.111.27'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.226.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.95%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '530.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.91%  '
$ws.Range("E7").Value = '  +1.26%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.219.63'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.609'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("E11").Value = '  -7.01%  '
$ws.Range("E12").Value = '  +3.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000255'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.13'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.739.43'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.52%  '
$ws.Range("E16").Value = '  -3.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.224.01'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.31'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '62.926.03'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.970'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '367.44'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.80'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.08'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.14'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.20'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '28.63'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '643.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.50'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.80%  '
$ws.Range("E34").Value = '  +3.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.106'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '57.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.93'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0719'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +17.97%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  +2.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.894.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.53'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +12.11%  '
$ws.Range("E46").Value = '  +3.94%  '
$ws.Range("E47").Value = '  +3.62%  '
$ws.Range("E48").Value = '  -1.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.33%  '
$ws.Range("E50").Value = '  +1.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '135.60'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.27%  '
